$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last (blank placeholder) row 12 entirely - rows shift up
$ws.Rows("12").Delete()

# Set row height for rows 10 and 11 to match the other data rows
$ws.Rows("10").RowHeight = 13.15
$ws.Rows("11").RowHeight = 13.15

# Fill in row 10 with new leaderboard data (D10 = 7/9/2025 as serial date)
$ws.Range("A10").Value = "VALLEY OFFICE PARK"
$ws.Range("B10").Value = "Steiner, Owen A"
$ws.Range("C10").Value = "003"
$ws.Range("D10").Value = 45847

# Fill in row 11 with new leaderboard data (D11 = 7/11/2025 as serial date)
$ws.Range("A11").Value = "BAKER AND LOVER INC"
$ws.Range("B11").Value = "Cina, Jonathan D"
$ws.Range("C11").Value = "023"
$ws.Range("D11").Value = 45849

# Copy date-formatted style from D6 onto D10/D11 so they render as dates
$ws.Range("D6").Copy()
$ws.Range("D10:D11").PasteSpecial(-4122) # xlPasteFormats

# Copy empty formatted style from E9 onto E10/E11 to match the other rows
$ws.Range("E9").Copy()
$ws.Range("E10:E11").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = $false

# Update selection to match final saved state
$ws.Range("B8").Select()
